$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted ahead of the existing
# "Vega Central Mapocho de Santiago - Zapallo italiano" records (row 302),
# pushing the rest of the rows (old 302-320) down by one (new 303-321).
$ws.Rows(302).Insert()

$ws.Cells.Item(302, 1).Value = 9
$ws.Cells.Item(302, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(302, 3).Value = "Metropolitana"
$ws.Cells.Item(302, 4).Value = 44610
$ws.Cells.Item(302, 5).Value = 13
$ws.Cells.Item(302, 6).Value = 100112032
$ws.Cells.Item(302, 7).Value = "Zapallo italiano"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 160
$ws.Cells.Item(302, 11).Value = 9000
$ws.Cells.Item(302, 12).Value = 10000
$ws.Cells.Item(302, 13).Value = 9500
$ws.Cells.Item(302, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(302, 15).Value = "Región Metropolitana"
$ws.Cells.Item(302, 16).Value = 190
$ws.Cells.Item(302, 17).Value = 50
$ws.Cells.Item(302, 18).Value = "Hortaliza"
